# Time-tracking Vinz - add new entry for 27.08. (REINFORCE algo, replay memory, Ubelix setup)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Fill in the new row of data (row 11 was blank before, between the last
# logged day on row 10 and the Total rows 16/17).
$ws.Range("A11").Value = "27.08."
$ws.Range("B11").Value = 555
$ws.Range("C11").Value = "Code REINFORCE algo, replay memory, setup Ubelix"

# Match the author's final selection / scroll position from the diff
# (scrolled so row 2 is the top visible row, with D7 as the active cell).
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D7").Select()
